$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 38 (existing rows 38-45 shift down to 39-46)
$ws.Rows.Item(38).Insert()

$newRow = $ws.Range("A38:D38")
$newRow.WrapText = $true
# Match the thin black border used by the surrounding rows (set Color before
# LineStyle so the engine reuses the existing border/style definitions
# instead of minting new ones).
$newRow.Borders.Color = 0
$newRow.Borders.LineStyle = 1

$ws.Range("A38").Value = "ENWIAM59"
$ws.Range("B38").Value = "OPQA-2924"
$ws.Range("C38").Value = "From ENW,verify that system is able to merge Activated STeAM account and Activated Facebook account and after merge verify STeAM TRUID is changed"
$ws.Range("D38").Value = "Y"

$ws.Rows.Item(38).RowHeight = 30

# Update the view: scroll position + active selection (best effort - the
# scroll position may not round-trip, but the selection does).
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E38").Select()
